$d = $word.ActiveDocument

function Get-ParagraphByPrefix($doc, $prefix) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# DESAFIO 11, 12 and 13 paragraphs get the same red (FF0000) run/paragraph-mark
# formatting already used by DESAFIO 5-10.
$p11 = Get-ParagraphByPrefix $d "DESAFIO 11:"
$p12 = Get-ParagraphByPrefix $d "DESAFIO 12:"
$p13 = Get-ParagraphByPrefix $d "DESAFIO 13:"

$p11.Range.Font.Color = 255
$p12.Range.Font.Color = 255
$p13.Range.Font.Color = 255

# The _GoBack bookmark used to sit at the end of the DESAFIO 11 paragraph;
# move it onto the trailing empty paragraph instead.
$b = $d.Bookmarks.Item("_GoBack")
$b.Delete()
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $last.Range
$r.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r)
